$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L10").Value = 367.8
$wsGrupo.Range("L18").Value = "3 de 16"

# Sheet: VENTA MENSUAL
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F10").Value = 367.8
$wsMensual.Range("F18").Value = 1619.87

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D2").Value = 6373.21
$wsCumplimiento.Range("E2").Value = -6373.21
$wsCumplimiento.Range("D4").Value = 20568.87
$wsCumplimiento.Range("E4").Value = -6845.53
$wsCumplimiento.Range("F4").Value = 1.49882390147005
